# attendance_sheet update: split combined percentage values, clarify
# "Not a part" label, and apply header/data formatting with conditional
# fill colors (yellow/orange/red) on the attendance columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Content fixes (do these first, before formatting, so the bulk
#    formatting calls below see a uniform range).
#    - row2: "133.33"/"100.00" were a mis-split combined value -> 66.67 / 50.00
#    - row4: "33.33" -> "16.67"
#    - "Not a part" -> "Not a part of class" (everywhere it appears)
# ---------------------------------------------------------------------
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "66.67"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "50.00"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "16.67"

$ws.Range("E3").Value = "Not a part of class"
$ws.Range("F3").Value = "Not a part of class"
$ws.Range("D4").Value = "Not a part of class"
$ws.Range("F4").Value = "Not a part of class"
$ws.Range("E5").Value = "Not a part of class"
$ws.Range("E6").Value = "Not a part of class"

# ---------------------------------------------------------------------
# 2. Header row (A1:F1): bold 12pt, centered, thin border around every
#    header cell except the last ("css A div") which gets no border.
# ---------------------------------------------------------------------
$header = $ws.Range("A1:F1")
$header.Font.Bold = $true
$header.Font.Size = 12
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108
$header.Borders.LineStyle = 1

$ws.Range("F1").Borders.LineStyle = -4142

$ws.Rows.Item(1).RowHeight = 20

# ---------------------------------------------------------------------
# 3. Data rows (A2:F6): Arial 10pt, centered, thin border as the
#    baseline look, then per-cell overrides below.
# ---------------------------------------------------------------------
$data = $ws.Range("A2:F6")
$data.Font.Name = "Arial"
$data.Font.Size = 10
$data.HorizontalAlignment = -4108
$data.VerticalAlignment = -4108
$data.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 4. Conditional fills on the attendance-percentage / status cells.
#    Yellow = good attendance, orange = borderline, red = low / not
#    part of class. A few column-F highlight cells sit without a
#    border (matches the original workbook's look).
# ---------------------------------------------------------------------
$ws.Range("E2").Interior.Color = 2681596   # FFFCEA28 yellow

$ws.Range("F2").Borders.LineStyle = -4142
$ws.Range("F2").Interior.Color = 36095     # FFFF8C00 orange

foreach ($c in @("D3", "E4", "D5", "D6")) {
    $ws.Range($c).Interior.Color = 3224055   # FFF73131 red
}

$ws.Range("F5").Borders.LineStyle = -4142
$ws.Range("F5").Interior.Color = 3224055   # FFF73131 red

$ws.Range("F6").Borders.LineStyle = -4142
$ws.Range("F6").Interior.Color = 3224055   # FFF73131 red

Write-Output "attendance sheet export formatting applied"
